# Update cryptos list values (price + 1h volume change) per source refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "28.294.36"
$ws.Cells.Item(2, 5).Value = "  -0.88%  "

$ws.Cells.Item(3, 4).Value = "1.550.50"
$ws.Cells.Item(3, 5).Value = "  -1.28%  "

$ws.Cells.Item(4, 5).Value = "  -0.04%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "209.38"
$ws.Cells.Item(5, 5).Value = "  -1.71%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.484"
$ws.Cells.Item(6, 5).Value = "  -1.54%  "

$ws.Cells.Item(7, 5).Value = "  -0.04%  "

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "23.57"
$ws.Cells.Item(8, 5).Value = "  -1.86%  "

$ws.Cells.Item(9, 5).Value = "  -2.01%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.0583"
$ws.Cells.Item(10, 5).Value = "  -1.30%  "

$ws.Cells.Item(11, 5).Value = "  +0.13%  "

$ws.Cells.Item(12, 4).Value = "1.771.17"
$ws.Cells.Item(12, 5).Value = "  -1.30%  "

$ws.Cells.Item(13, 4).Value = "1.553.34"
$ws.Cells.Item(13, 5).Value = "  -1.10%  "

$ws.Cells.Item(14, 4).Value = "28.284.71"
$ws.Cells.Item(14, 5).Value = "  -0.92%  "

$ws.Cells.Item(15, 5).Value = "  -1.52%  "

$ws.Cells.Item(16, 5).Value = "  -2.31%  "

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "60.53"
$ws.Cells.Item(17, 5).Value = "  -2.79%  "

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "226.92"
$ws.Cells.Item(18, 5).Value = "  -1.28%  "

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "7.31"
$ws.Cells.Item(19, 5).Value = "  -0.71%  "

$ws.Cells.Item(20, 4).Value = "0.0₃0673"

$ws.Cells.Item(21, 5).Value = "  -0.03%  "

$ws.Cells.Item(22, 5).Value = "  +0.89%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "8.83"
$ws.Cells.Item(23, 5).Value = "  -3.15%  "

$ws.Cells.Item(24, 5).Value = "  -5.67%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "148.75"
$ws.Cells.Item(25, 5).Value = "  -2.02%  "

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "14.73"
$ws.Cells.Item(26, 5).Value = "  -1.99%  "

$ws.Cells.Item(27, 5).Value = "  -0.41%  "

$ws.Cells.Item(28, 5).Value = "  -0.07%  "

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "6.23"
$ws.Cells.Item(29, 5).Value = "  -3.16%  "

$ws.Cells.Item(30, 5).Value = "  -4.08%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "1.06"
$ws.Cells.Item(31, 5).Value = "  -4.52%  "

$ws.Cells.Item(32, 5).Value = "  -1.16%  "

$ws.Cells.Item(33, 5).Value = "  -1.83%  "

$ws.Cells.Item(34, 4).Value = "1.382.00"
$ws.Cells.Item(34, 5).Value = "  -0.88%  "

$ws.Cells.Item(35, 5).Value = "  +0.40%  "

$ws.Cells.Item(36, 5).Value = "  -3.44%  "

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "2.34"
$ws.Cells.Item(37, 5).Value = "  -1.27%  "

$ws.Cells.Item(38, 5).Value = "  -1.40%  "

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.0160"
$ws.Cells.Item(39, 5).Value = "  -2.93%  "

$ws.Cells.Item(40, 2).Value = "PaxDollar"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "1.00"
$ws.Cells.Item(40, 5).Value = "  -0.03%  "

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "1.90"
$ws.Cells.Item(41, 5).Value = "  +0.52%  "

$ws.Cells.Item(42, 2).Value = "ImmutableX"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.507"
$ws.Cells.Item(42, 5).Value = "  -3.19%  "

$ws.Cells.Item(43, 5).Value = "  -1.86%  "

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.0465"
$ws.Cells.Item(44, 5).Value = "  -1.16%  "

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "5.38"
$ws.Cells.Item(45, 5).Value = "  -2.19%  "

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "61.73"
$ws.Cells.Item(46, 5).Value = "  -1.86%  "

$ws.Cells.Item(47, 2).Value = "WEMIXToken"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.907"
$ws.Cells.Item(47, 5).Value = "  -6.36%  "

$ws.Cells.Item(48, 2).Value = "RocketPoolETH"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Cells.Item(48, 4).Value = "1.684.56"
$ws.Cells.Item(48, 5).Value = "  -1.31%  "

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "85.23"
$ws.Cells.Item(49, 5).Value = "  -1.24%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "42.35"
$ws.Cells.Item(50, 5).Value = "  +7.12%  "

$ws.Cells.Item(51, 5).Value = "  +0.15%  "
